# Update the AQMOS waterontharder capacity calculation workbook:
# switch the selected softener model from LESS-10 to LESS-15 and
# move the active selection to D25 on the Dashboard sheet.

$wb = $excel.ActiveWorkbook

$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Activate()

# The model picker (B3) is a data-validation dropdown sourced from
# Data!$B$2:$B$8. Setting the value directly mimics picking "LESS-15"
# from that list; every downstream formula (Data!F2/G2, the K/L lookup
# table, D12/D13, ...) recalculates automatically from this single input.
$dashboard.Range("B3").Value = "LESS-15"

# Move the selection the way the author last left it in the sheet.
$dashboard.Range("D25").Select()

$excel.CalculateFullRebuild()
